# Fix cellular filtering configuration on the "Main Info" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Info")

# APN (backup): internet.odjosky.com -> internet.odferky.com
$ws.Range("B26").Value = "internet.odferky.com"

# Design: FLOW -> BASE
$ws.Range("B7").Value = "BASE"

# Region: EMEA -> NAM
$ws.Range("B2").Value = "NAM"

# Converged router: FALSE -> TRUE
$ws.Range("B8").Value = $true

# Migration from MPLS: "True - Production router" -> FALSE
$ws.Range("B9").Value = $false

# 4G+Cellular (backup): TRUE -> FALSE
$ws.Range("B25").Value = $false

# Update the active selection to D15
$ws.Range("D15").Select()
